$wb = $excel.ActiveWorkbook

# --- Update "Score" sheet ---
$scoreSheet = $wb.Worksheets.Item("Score")

# Update the shared ratio/description strings (B2, C2)
$scoreSheet.Range("B2").Value = "0.02:0.01:0.03:0.03:0.91"
$scoreSheet.Range("C2").Value = "0.02T, 0.01RR, 0.03Env, 0.03Econ, 0.91S"

# Update the score values for each alternative (D2, E2, F2)
$scoreSheet.Range("D2").Value = 0.3813762478735621
$scoreSheet.Range("E2").Value = 0.9783077628436179
$scoreSheet.Range("F2").Value = 0.009442416060179716

# --- Update "Rank" sheet ---
$rankSheet = $wb.Worksheets.Item("Rank")

# Update the shared ratio/description strings (B2, C2) to stay consistent
$rankSheet.Range("B2").Value = "0.02:0.01:0.03:0.03:0.91"
$rankSheet.Range("C2").Value = "0.02T, 0.01RR, 0.03Env, 0.03Econ, 0.91S"

# Update the rank values (swap D2 and F2)
$rankSheet.Range("D2").Value = 2
$rankSheet.Range("F2").Value = 3
